# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (positioned between "2021-Q4" and "总计")
# with the detailed per-fund holding breakdown, and inserts a new leading
# row into the "总计" (totals) sheet summarising the 2022-Q1 quarter while
# keeping the pre-existing 2021-Q4 summary row beneath it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Header row - same layout as the "2021-Q4" sheet. Copy the header cell's
# format (bold font + border + centered alignment, style index 2 in the
# original workbook) onto the new header cells so they match visually.
$q4.Cells.Item(1, 2).Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Data row. Reuse the "index" cell's style (A2 on the "2021-Q4" sheet) for
# the leading A2 cell on the new sheet too.
$q4.Cells.Item(2, 1).Copy()
$q1.Cells.Item(2, 1).PasteSpecial(-4122)
$q1.Cells.Item(2, 1).Value = 0

# The fund code has a significant leading zero, so it must stay text too.
# (Temporarily force Text format while assigning, then drop back to the
# Normal style so the cell's format index matches the plain unstyled cells
# around it - otherwise Excel auto-converts the numeric-looking string and
# loses the leading zero.)
$q1.Cells.Item(2, 2).NumberFormat = "@"
$q1.Cells.Item(2, 2).Value = "006105"
$q1.Cells.Item(2, 2).Style = "Normal"
$q1.Cells.Item(2, 3).Value = "泰达宏利印度机会股票（QDII）"

# These columns hold formatted decimal strings (e.g. "0.60") where the
# trailing/leading zeros are significant, so force Text format before
# assigning - otherwise Excel auto-converts them to plain numbers.
$q1.Cells.Item(2, 4).NumberFormat = "@"
$q1.Cells.Item(2, 4).Value = "0.60"
$q1.Cells.Item(2, 4).Style = "Normal"
$q1.Cells.Item(2, 5).NumberFormat = "@"
$q1.Cells.Item(2, 5).Value = "87.31"
$q1.Cells.Item(2, 5).Style = "Normal"
$q1.Cells.Item(2, 6).NumberFormat = "@"
$q1.Cells.Item(2, 6).Value = "4.20"
$q1.Cells.Item(2, 6).Style = "Normal"
$q1.Cells.Item(2, 7).NumberFormat = "@"
$q1.Cells.Item(2, 7).Value = "0.0252"
$q1.Cells.Item(2, 7).Style = "Normal"

$q1.Cells.Item(2, 8).Value = 5

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for 2022-Q1 above the
#    existing 2021-Q4 summary row (which shifts down to row 3).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Rows.Item(2).ClearFormats()

# Restore the "index" column's style (bold/bordered, same as the header and
# as the row below it) on the new A2 cell - ClearFormats() above wiped it.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.03

# The row-insert keeps the old row's values but its index cell (A3) needs
# bumping from 0 to 1 now that it is the second data row.
$total.Cells.Item(3, 1).Value = 1
